$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47; this shifts existing rows 47-65 down to 48-66
$ws.Rows.Item(47).Insert()

# Match the date formatting used by the other rows in column D before setting the value
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat

# Populate the new row 47 with the new weekly data entry
$ws.Cells.Item(47, 1).Value = 9
$ws.Cells.Item(47, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(47, 3).Value = "Metropolitana"
$ws.Cells.Item(47, 4).Value = 44463
$ws.Cells.Item(47, 5).Value = 13
$ws.Cells.Item(47, 6).Value = 100112005
$ws.Cells.Item(47, 7).Value = "Puerro"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 160
$ws.Cells.Item(47, 11).Value = 7500
$ws.Cells.Item(47, 12).Value = 8000
$ws.Cells.Item(47, 13).Value = 7750
$ws.Cells.Item(47, 14).Value = "`$/paquete 20 unidades"
$ws.Cells.Item(47, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(47, 16).Value = 388
$ws.Cells.Item(47, 17).Value = 20
$ws.Cells.Item(47, 18).Value = "Hortaliza"
